# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the status moves from
# "In Translation" to "Ready for handoff" and the associated timestamps are
# bumped forward by under a minute (the new Xliff handoff generation run).
#
# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
# zh-cn / de-de sheets: Status column (C2) and Latest Handoff Datetime (H2).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 14:48:47"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 14:48:42"

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 14:48:47"

# The Status columns now hold the longer "Ready for handoff" string, so the
# report widens those columns to fit it (matches the regenerated report's
# layout): Overview columns E & F, and the "Status" column (C) on the
# zh-cn / de-de sheets.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
